$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles, number formats) from the last existing data row (62)
# into the new row (63) so the new row matches the look of prior rows.
$ws.Range("A62:R62").Copy($ws.Range("A63:R63"))

# Now populate the new row with the actual data values described in the diff.
$ws.Range("A63").Value = 4
$ws.Range("B63").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C63").Value = "Los Lagos"
$ws.Range("D63").Value = 44628
$ws.Range("E63").Value = 10
$ws.Range("F63").Value = 100112031
$ws.Range("G63").Value = "Poroto verde"
$ws.Range("H63").Value = "Magnum"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 50
$ws.Range("K63").Value = 38000
$ws.Range("L63").Value = 38000
$ws.Range("M63").Value = 38000
$ws.Range("N63").Value = "$/saco 25 kilos"
$ws.Range("O63").Value = "Región Metropolitana"
$ws.Range("P63").Value = 1520
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
